# Apply the OOXML diff: rename layer names in column B (rows 8-21) to
# their ONNX-style node output names, and update the detection-head output
# channel counts (column O/L and the time column S) plus the grand total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B layer-name renames (rows 8..21) ---
$ws.Range("B8").Value  = "/module_list.0/conv_0/Conv_output_0"
$ws.Range("B9").Value  = "/module_list.1/maxpool_1/MaxPool_output_0"
$ws.Range("B10").Value = "/module_list.2/conv_2/Conv_output_0"
$ws.Range("B11").Value = "/module_list.3/maxpool_3/MaxPool_output_0"
$ws.Range("B12").Value = "/module_list.4/conv_4/Conv_output_0"
$ws.Range("B13").Value = "/module_list.5/maxpool_5/MaxPool_output_0"
$ws.Range("B14").Value = "/module_list.6/conv_6/Conv_output_0"
$ws.Range("B15").Value = "/module_list.7/maxpool_7/MaxPool_output_0"
$ws.Range("B16").Value = "/module_list.8/conv_8/Conv_output_0"
$ws.Range("B17").Value = "/module_list.9/maxpool_9/MaxPool_output_0"
$ws.Range("B18").Value = "/module_list.10/conv_10/Conv_output_0"
$ws.Range("B19").Value = "/module_list.11/maxpool_11/MaxPool_output_0"
$ws.Range("B20").Value = "/module_list.12/conv_12/Conv_output_0"
$ws.Range("B21").Value = "/module_list.13/conv_13/Conv_output_0"

# --- Row 22 (output1 / Conv): output channels 125 -> 30, time 136 -> 126 ---
$ws.Range("O22").Value = 30
$ws.Range("S22").Value = 126

# --- Row 23 (post_0_0_transpose): in/out channels 125 -> 30, time 20 -> 5 ---
$ws.Range("L23").Value = 30
$ws.Range("O23").Value = 30
$ws.Range("S23").Value = 5

# --- Row 24 (post_0_1_cast_fp16_fp32): in/out channels 125 -> 30, time 17 -> 4 ---
$ws.Range("L24").Value = 30
$ws.Range("O24").Value = 30
$ws.Range("S24").Value = 4

# --- Row 25 (Total): time 33846 -> 33808 ---
$ws.Range("S25").Value = 33808
